$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 1625.7142
$ws.Range("I40").Value = 1731.6666
$ws.Range("K40").Value = 1731.6666
$ws.Range("M40").Value = -1556.6666

# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 3020
$ws.Range("J58").Value = 10000
$ws.Range("L58").Value = 30000
$ws.Range("N58").Value = -30300

# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 2575.8635
$ws.Range("I100").Value = 2575.8635
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2575.8635
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2034.8635
$ws.Range("N100").Value = ""

# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 383.375
$ws.Range("I107").Value = 352.42856
$ws.Range("K107").Value = 352.42856
$ws.Range("M107").Value = 1567.57144

# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 800.375
$ws.Range("I125").Value = 832
$ws.Range("K125").Value = 7488
$ws.Range("M125").Value = -5028

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2327.1
$ws.Range("I137").Value = 1633.75
$ws.Range("J137").Value = 2789.3333
$ws.Range("K137").Value = 4901.25
$ws.Range("L137").Value = 8367.999899999999
$ws.Range("M137").Value = -2351.25
$ws.Range("N137").Value = -13467.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 2565.2
$ws.Range("I32").Value = 2011.421
$ws.Range("K32").Value = 2011.421
$ws.Range("M32").Value = -1724.421

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 2072.875
$ws.Range("I61").Value = 1797.5714
$ws.Range("K61").Value = 1797.5714
$ws.Range("M61").Value = -1585.5714

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 4072.182
$ws.Range("I74").Value = 3310
$ws.Range("J74").Value = 4986.8
$ws.Range("K74").Value = 3310
$ws.Range("L74").Value = 4986.8
$ws.Range("M74").Value = -2436
$ws.Range("N74").Value = -6734.8

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 4072.182
$ws.Range("I77").Value = 3310
$ws.Range("J77").Value = 4986.8
$ws.Range("K77").Value = 16550
$ws.Range("L77").Value = 24934
$ws.Range("M77").Value = -12182
$ws.Range("N77").Value = -33670

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 3224.2
$ws.Range("I102").Value = 3155.25
$ws.Range("K102").Value = 3155.25
$ws.Range("M102").Value = -1533.25

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 7232.871
$ws.Range("I132").Value = 6229.476
$ws.Range("K132").Value = 18688.428
$ws.Range("M132").Value = -16158.428

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 2072.875
$ws.Range("I136").Value = 1797.5714
$ws.Range("K136").Value = 5392.7142
$ws.Range("M136").Value = -2842.7142

$ws = $wb.Worksheets.Item("BSM")
# Row 13 (Leve Item ID 27127)
$ws.Range("H13").Value = 40000
$ws.Range("J13").Value = 40000
$ws.Range("L13").Value = 40000
$ws.Range("N13").Value = -40336

# Row 96 (Leve Item ID 19525)
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").Value = ""

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 3364.8
$ws.Range("I107").Value = 2724.6667
$ws.Range("K107").Value = 2724.6667
$ws.Range("M107").Value = -804.6667000000002

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 112.77778
$ws.Range("I7").Value = 69.375
$ws.Range("J7").Value = 460
$ws.Range("K7").Value = 69.375
$ws.Range("L7").Value = 460
$ws.Range("M7").Value = 43.625
$ws.Range("N7").Value = -686

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 3220.1667
$ws.Range("I58").Value = 2222
$ws.Range("K58").Value = 2222
$ws.Range("M58").Value = -2019

# Row 59 (Leve Item ID 1942)
$ws.Range("H59").Value = 61333.25
$ws.Range("J59").Value = 75076.336
$ws.Range("L59").Value = 75076.336
$ws.Range("N59").Value = -77366.336

# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 13076.111
$ws.Range("I86").Value = 6812.2856
$ws.Range("K86").Value = 6812.2856
$ws.Range("M86").Value = -5689.2856

# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 13076.111
$ws.Range("I89").Value = 6812.2856
$ws.Range("K89").Value = 34061.428
$ws.Range("M89").Value = -28445.428

# Row 116 (Leve Item ID 26117)
$ws.Range("H116").Value = 56362.816
$ws.Range("J116").Value = 56362.816
$ws.Range("L116").Value = 56362.816
$ws.Range("N116").Value = -65540.81599999999

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1088.909
$ws.Range("I134").Value = 1088.909
$ws.Range("K134").Value = 3266.727
$ws.Range("M134").Value = -731.7270000000003

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 3220.1667
$ws.Range("I136").Value = 2222
$ws.Range("K136").Value = 6666
$ws.Range("M136").Value = -4116

# Row 141 (Leve Item ID 43345)
$ws.Range("H141").Value = 707516.3
$ws.Range("J141").Value = 707516.3
$ws.Range("L141").Value = 707516.3
$ws.Range("N141").Value = -717876.3

$ws = $wb.Worksheets.Item("CUL")
# Row 92 (Leve Item ID 19841)
$ws.Range("H92").Value = 362.125
$ws.Range("I92").Value = 362.125
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1086.375
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 161.625
$ws.Range("N92").Value = ""

# Row 98 (Leve Item ID 19843)
$ws.Range("H98").Value = 495
$ws.Range("J98").Value = 495
$ws.Range("L98").Value = 1485
$ws.Range("N98").Value = -4481

# Row 121 (Leve Item ID 27878)
$ws.Range("H121").Value = 3452.7144
$ws.Range("J121").Value = 3452.7144
$ws.Range("L121").Value = 10358.1432
$ws.Range("N121").Value = -12978.1432

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2072
$ws.Range("I132").Value = 2082.2856
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 6246.8568
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3716.8568
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 4089.6553
$ws.Range("I46").Value = 4089.6553
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 4089.6553
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -3901.6553
$ws.Range("N46").Value = ""

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 1280.7858
$ws.Range("I55").Value = 409.25
$ws.Range("K55").Value = 409.25
$ws.Range("M55").Value = -236.25

# Row 116 (Leve Item ID 26133)
$ws.Range("H116").Value = 35000
$ws.Range("J116").Value = 35000
$ws.Range("L116").Value = 35000
$ws.Range("N116").Value = -44178

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 4387.778
$ws.Range("I132").Value = 3927.1428
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 11781.4284
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -9251.428400000001
$ws.Range("N132").Value = -23060

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 2929.3333
$ws.Range("I136").Value = 2329
$ws.Range("J136").Value = 4130
$ws.Range("K136").Value = 6987
$ws.Range("L136").Value = 12390
$ws.Range("M136").Value = -4437
$ws.Range("N136").Value = -17490

$ws = $wb.Worksheets.Item("WVR")
# Row 4 (Leve Item ID 2996)
$ws.Range("H4").Value = 8600
$ws.Range("J4").Value = 8600
$ws.Range("L4").Value = 8600
$ws.Range("N4").Value = -8826

# Row 54 (Leve Item ID 3413)
$ws.Range("H54").Value = 14947.35
$ws.Range("J54").Value = 29589.4
$ws.Range("L54").Value = 29589.4
$ws.Range("N54").Value = -30629.4

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 833.64703
$ws.Range("I107").Value = 825.3333
$ws.Range("K107").Value = 2475.9999
$ws.Range("M107").Value = -555.9998999999998

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2512.5715
$ws.Range("I132").Value = 2765.3333
$ws.Range("K132").Value = 8295.999899999999
$ws.Range("M132").Value = -5765.999899999999

# Row 140 (Leve Item ID 42506)
$ws.Range("H140").Value = 66944
$ws.Range("J140").Value = 66944
$ws.Range("L140").Value = 66944
$ws.Range("N140").Value = -77304
